# Convert EMU -> points for the COM layer. A tiny epsilon is added because
# the host stores shape Left/Top/Width/Height as (lossy) floating point
# points; without the nudge values land 1 EMU short after the point->EMU
# round trip.
function EMU($v) {
    return ($v / 12700) + 0.00001
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# "Scikit-learn" text placeholder (id 3) - give it an explicit position/size.
$shScikitLabel = $s.Shapes.Item(2)
$shScikitLabel.Left   = EMU(827088)
$shScikitLabel.Top    = EMU(1826039)
$shScikitLabel.Width  = EMU(5157787)
$shScikitLabel.Height = EMU(823912)

# Scikit-learn content placeholder (id 4).
$shScikitBody = $s.Shapes.Item(3)
$shScikitBody.Left   = EMU(827088)
$shScikitBody.Top    = EMU(2649951)
$shScikitBody.Width  = EMU(5157787)
$shScikitBody.Height = EMU(3684588)

# "statsmodels" text placeholder (id 5).
$shStatsLabel = $s.Shapes.Item(4)
$shStatsLabel.Left   = EMU(6159500)
$shStatsLabel.Top    = EMU(1826039)
$shStatsLabel.Width  = EMU(5183188)
$shStatsLabel.Height = EMU(823912)

# statsmodels content placeholder (id 6).
$shStatsBody = $s.Shapes.Item(5)
$shStatsBody.Left   = EMU(6159500)
$shStatsBody.Top    = EMU(2649951)
$shStatsBody.Width  = EMU(5183188)
$shStatsBody.Height = EMU(3684588)

# New textbox (id 7) describing the test setup, added below the title.
$newBoxLeft   = EMU(827088)
$newBoxTop    = EMU(1435893)
$newBoxWidth  = EMU(10515600)
$newBoxHeight = EMU(1325563)
$newBox = $s.Shapes.AddTextbox(1, $newBoxLeft, $newBoxTop, $newBoxWidth, $newBoxHeight)
$newBox.Name = "Title 1"
$newBox.TextFrame.TextRange.Text = "Test: 1000 samples, 20 features (independent variables), 2 groups"
$newBox.TextFrame.TextRange.Font.Size = 24

$newBox.TextFrame2.Orientation = 1
$newBox.TextFrame.MarginLeft   = EMU(91440)
$newBox.TextFrame.MarginTop    = EMU(45720)
$newBox.TextFrame.MarginRight  = EMU(91440)
$newBox.TextFrame.MarginBottom = EMU(45720)
$newBox.TextFrame.VerticalAnchor = 1
$newBox.TextFrame.AutoSize = 2
